$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the (only) worksheet from the generic "1" to the municipality name.
$ws.Name = "ყაზბეგი"

# Drop the "(მოსახლეობის აღწერის შედეგებით)" sub-heading row entirely -
# it is no longer part of the exported table.
$ws.Rows(2).Delete()

# The table used to report area for three census years (1989 / 2002 / 2014);
# keep only the most recent (2014) column and drop the other two.
$ws.Columns("B:C").Delete()

# Match the saved selection of the authored workbook.
$null = $ws.Range("A2").Select()
